$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark from its old location
#    (it sat just before the run containing "18.72" in the
#    "Overall Precision:" paragraph of the "Last 92 Text Files"
#    section).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Turn the first empty paragraph that follows the "Future
#    Improvements" subtitle into a "No Spacing" paragraph containing
#    the reflection text, and re-create the "_GoBack" bookmark as a
#    collapsed range right after that text (inside the paragraph,
#    before the paragraph mark) - matching what Word does when the
#    cursor is left there on save.
# ------------------------------------------------------------------
$target = $null
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $cur = $d.Paragraphs.Item($i)
    if ($cur.Range.Text.TrimEnd([char]13, [char]7) -eq "Future Improvements") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

$target.Style = "No Spacing"

$reflection = "I found myself spending more time working on entity tagging then ontology construction. If I were to re-do this assignments I" + [char]0x2019 + "d ensure I split my time more equally between the two tasks."

# A small amount of throwaway padding is appended after the real text
# before inserting, then stripped again once the bookmark has been
# anchored.  Placing a *collapsed* bookmark range right at the very
# tail of the document body confuses the host's range resolution, so
# the padding keeps the anchor point comfortably away from the
# document end while the bookmark is created.
$padding = "XXX"

$r = $target.Range
$r.InsertBefore($reflection + $padding)

$anchorPos = $r.End - 1 - $padding.Length
$bmRange = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padRange = $d.Range($anchorPos, $anchorPos + $padding.Length)
$padRange.Delete()
